$wb = $excel.ActiveWorkbook

# Sheet "5_" (the sampling-vocabulary matching exercise) is the sheet whose
# numeric answer cells get converted to quote-prefixed text values, and which
# becomes the new active tab/sheet.
$ws = $wb.Worksheets.Item("5_")

# --- Pre-create the two new cell styles in the exact index order Excel used ---
# Style index 6 = (fillId 3, wrapText) + quotePrefix  -- derived from C2's style
# Style index 7 = (fillId 3, no wrap)   + quotePrefix  -- derived from C6's style
# We provoke their creation (in this order) via two throwaway cells well
# outside the sheet's used range, then remove those cells entirely so they
# leave no trace in the saved worksheet XML.
$ws.Range("C2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z100").Value = "''seed1'"

$ws.Range("C6").Copy()
$ws.Range("Z101").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z101").Value = "''seed2'"

$ws.Range("Z100:Z101").Clear()

# --- Replace the numeric sample values with quote-prefixed text equivalents ---
# Written in the same order the shared-string table ends up in (C7, C8, then
# C2..C6), so the new <si> entries land at the expected indices.
$ws.Range("C7").Value = "''4000'"
$ws.Range("C8").Value = "''0.005'"
$ws.Range("C2").Value = "''1000'"
$ws.Range("C3").Value = "''200'"
$ws.Range("C4").Value = "''4'"
$ws.Range("C5").Value = "''800'"
$ws.Range("C6").Value = "''0.001'"

# --- Move the selected/active tab from "10_" to "5_" ---
$ws.Activate()
[void]$ws.Range("C11").Select()

Write-Output "edit complete"
